# Incomes sheet: a new "Salary" income entry (10000, dated 2025-07-24) is
# added as row 4. The rows that used to follow ("Trading" / "Business")
# shift down by one (old row4->5, old row5->6); the used range grows from
# A1:C5 to A1:C6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 4, pushing "Trading" (row4) and "Business"
# (row5) down by one row each.
$ws.Rows(4).Insert()

# Populate the newly inserted row 4 with the new income entry.
$ws.Range("A4").Value = "Salary"
$ws.Range("B4").Value = 10000
$ws.Range("C4").Value = 45862.22928240741

# Match the date formatting used by the other Date-column cells (style is
# copied from C3 so no new/duplicate number format gets created).
$ws.Range("C3").Copy()
$ws.Range("C4").PasteSpecial(-4122)
$excel.CutCopyMode = 0
